# "Lista De Pascale, nuovi grafici"
# Adds a new "DE PASCALE" list row to the COALIZIONE/LISTA/COLORE table on
# Foglio1, renames "ALLEANZA VERDI E SINISTRA" -> "ALLEANZA VERDI SINISTRA"
# for that row, and reclassifies the former M5S/CENTRO coalition rows as
# SINISTRA (the coalition supporting De Pascale). Also restores Foglio1 as
# the active sheet/selection.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert a new row above the current row 3 ("ALLEANZA VERDI E SINISTRA"),
# shifting it (and everything below) down by one.
$ws1.Rows.Item(3).Insert()

# Populate the new row 3 with the De Pascale list entry.
$ws1.Range("A3").Value = "SINISTRA"
$ws1.Range("B3").Value = "DE PASCALE"
$ws1.Range("C3").Value = "#DDA530"

# Old row 3 ("ALLEANZA VERDI E SINISTRA") is now row 4: drop the "E".
$ws1.Range("B4").Value = "ALLEANZA VERDI SINISTRA"

# Old rows 4-6 (M5S / CENTRO / CENTRO coalitions) are now rows 5-7: they all
# move into the SINISTRA coalition alongside the PD and the other allies.
$ws1.Range("A5").Value = "SINISTRA"
$ws1.Range("A6").Value = "SINISTRA"
$ws1.Range("A7").Value = "SINISTRA"

# Restore Foglio1 as the active/selected sheet with B5 selected.
$null = $ws1.Select()
$null = $ws1.Range("B5").Select()
